# Generate Report for Handoff
# Adds two new localization entries (646596ff-... and 68bcec46-...) to the
# Overview / zh-cn / de-de report sheets, mirroring the existing rows.

$wb = $excel.ActiveWorkbook

$dateFormat = "yyyy-mm-dd HH:mm:ss"

# ---------------------------------------------------------------------------
# Sheet "Overview" (sheet1 / table3)
# ---------------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")
$loOverview = $wsOverview.ListObjects.Item(1)
$loOverview.ListRows.Add() | Out-Null
$loOverview.ListRows.Add() | Out-Null

$wsOverview.Range("A4").Value = "646596ff-78c4-4b47-b2ef-2067e79379b0.md"
$wsOverview.Range("B4").Value = "e2e\646596ff-78c4-4b47-b2ef-2067e79379b0.md"
$wsOverview.Range("C4").Value = ".md"
$wsOverview.Range("D4").Value = "'"
$wsOverview.Range("E4").Value = "Ready for handoff"
$wsOverview.Range("F4").Value = "Ready for handoff"
$wsOverview.Range("G4").Value = "2016-08-20 08:49:02"
$wsOverview.Range("G4").NumberFormat = $dateFormat

$wsOverview.Range("A5").Value = "68bcec46-6cb0-401a-b119-925d1709bc41.md"
$wsOverview.Range("B5").Value = "e2e\68bcec46-6cb0-401a-b119-925d1709bc41.md"
$wsOverview.Range("C5").Value = ".md"
$wsOverview.Range("D5").Value = "'"
$wsOverview.Range("E5").Value = "Ready for handoff"
$wsOverview.Range("F5").Value = "Ready for handoff"
$wsOverview.Range("G5").Value = "2016-08-20 08:49:02"
$wsOverview.Range("G5").NumberFormat = $dateFormat

$wsOverview.Hyperlinks.Add($wsOverview.Range("B4"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/3267b74449032ebf311bf224864d8f9f9bc063be/e2e/646596ff-78c4-4b47-b2ef-2067e79379b0.md", "", "", "e2e\646596ff-78c4-4b47-b2ef-2067e79379b0.md") | Out-Null
$wsOverview.Hyperlinks.Add($wsOverview.Range("B5"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/3267b74449032ebf311bf224864d8f9f9bc063be/e2e/68bcec46-6cb0-401a-b119-925d1709bc41.md", "", "", "e2e\68bcec46-6cb0-401a-b119-925d1709bc41.md") | Out-Null

# ---------------------------------------------------------------------------
# Sheet "zh-cn" (sheet2 / table1)
# ---------------------------------------------------------------------------
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$loZhCn = $wsZhCn.ListObjects.Item(1)
$loZhCn.ListRows.Add() | Out-Null
$loZhCn.ListRows.Add() | Out-Null

$wsZhCn.Range("A4").Value = "646596ff-78c4-4b47-b2ef-2067e79379b0.md"
$wsZhCn.Range("B4").Value = ".md"
$wsZhCn.Range("C4").Value = "Ready for handoff"
$wsZhCn.Range("D4").Value = "e2e"
$wsZhCn.Range("E4").Value = "ht"
$wsZhCn.Range("F4").Value = "'False"
$wsZhCn.Range("G4").Value = "646596ff-78c4-4b47-b2ef-2067e79379b0.f8dee69f46c85309e3a254a3182fcb4c0cdbd8ba.zh-cn.xlf"
$wsZhCn.Range("H4").Value = "2016-08-20 08:48:56"
$wsZhCn.Range("H4").NumberFormat = $dateFormat
$wsZhCn.Range("I4").Value = "'"
$wsZhCn.Range("J4").Value = "'"
$wsZhCn.Range("K4").Value = "0001-01-01 00:00:00"
$wsZhCn.Range("K4").NumberFormat = $dateFormat
$wsZhCn.Range("L4").Value = "'"
$wsZhCn.Range("M4").Value = "'True"
$wsZhCn.Range("N4").Value = "'"
$wsZhCn.Range("O4").Value = "'False"
$wsZhCn.Range("P4").Value = "'"

$wsZhCn.Range("A5").Value = "68bcec46-6cb0-401a-b119-925d1709bc41.md"
$wsZhCn.Range("B5").Value = ".md"
$wsZhCn.Range("C5").Value = "Ready for handoff"
$wsZhCn.Range("D5").Value = "e2e"
$wsZhCn.Range("E5").Value = "ht"
$wsZhCn.Range("F5").Value = "'False"
$wsZhCn.Range("G5").Value = "68bcec46-6cb0-401a-b119-925d1709bc41.a809f9eebb6c99ea08094d75a9e956913b9d5c46.zh-cn.xlf"
$wsZhCn.Range("H5").Value = "2016-08-20 08:48:56"
$wsZhCn.Range("H5").NumberFormat = $dateFormat
$wsZhCn.Range("I5").Value = "'"
$wsZhCn.Range("J5").Value = "'"
$wsZhCn.Range("K5").Value = "0001-01-01 00:00:00"
$wsZhCn.Range("K5").NumberFormat = $dateFormat
$wsZhCn.Range("L5").Value = "'"
$wsZhCn.Range("M5").Value = "'True"
$wsZhCn.Range("N5").Value = "'"
$wsZhCn.Range("O5").Value = "'False"
$wsZhCn.Range("P5").Value = "'"

$wsZhCn.Hyperlinks.Add($wsZhCn.Range("A4"), "https://github.com/OpenLocalizationTestOrg/ol-test0-zhcn/blob/23458073678df38eb3f7fa6b421d5480117f608b/e2e/646596ff-78c4-4b47-b2ef-2067e79379b0.md", "", "", "646596ff-78c4-4b47-b2ef-2067e79379b0.md") | Out-Null
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("A5"), "https://github.com/OpenLocalizationTestOrg/ol-test0-zhcn/blob/23458073678df38eb3f7fa6b421d5480117f608b/e2e/68bcec46-6cb0-401a-b119-925d1709bc41.md", "", "", "68bcec46-6cb0-401a-b119-925d1709bc41.md") | Out-Null

# ---------------------------------------------------------------------------
# Sheet "de-de" (sheet3 / table2)
# ---------------------------------------------------------------------------
$wsDeDe = $wb.Worksheets.Item("de-de")
$loDeDe = $wsDeDe.ListObjects.Item(1)
$loDeDe.ListRows.Add() | Out-Null
$loDeDe.ListRows.Add() | Out-Null

$wsDeDe.Range("A4").Value = "646596ff-78c4-4b47-b2ef-2067e79379b0.md"
$wsDeDe.Range("B4").Value = ".md"
$wsDeDe.Range("C4").Value = "Ready for handoff"
$wsDeDe.Range("D4").Value = "e2e"
$wsDeDe.Range("E4").Value = "ht"
$wsDeDe.Range("F4").Value = "'False"
$wsDeDe.Range("G4").Value = "646596ff-78c4-4b47-b2ef-2067e79379b0.f8dee69f46c85309e3a254a3182fcb4c0cdbd8ba.de-de.xlf"
$wsDeDe.Range("H4").Value = "2016-08-20 08:49:02"
$wsDeDe.Range("H4").NumberFormat = $dateFormat
$wsDeDe.Range("I4").Value = "'"
$wsDeDe.Range("J4").Value = "'"
$wsDeDe.Range("K4").Value = "0001-01-01 00:00:00"
$wsDeDe.Range("K4").NumberFormat = $dateFormat
$wsDeDe.Range("L4").Value = "'"
$wsDeDe.Range("M4").Value = "'True"
$wsDeDe.Range("N4").Value = "'"
$wsDeDe.Range("O4").Value = "'False"
$wsDeDe.Range("P4").Value = "'"

$wsDeDe.Range("A5").Value = "68bcec46-6cb0-401a-b119-925d1709bc41.md"
$wsDeDe.Range("B5").Value = ".md"
$wsDeDe.Range("C5").Value = "Ready for handoff"
$wsDeDe.Range("D5").Value = "e2e"
$wsDeDe.Range("E5").Value = "ht"
$wsDeDe.Range("F5").Value = "'False"
$wsDeDe.Range("G5").Value = "68bcec46-6cb0-401a-b119-925d1709bc41.a809f9eebb6c99ea08094d75a9e956913b9d5c46.de-de.xlf"
$wsDeDe.Range("H5").Value = "2016-08-20 08:49:02"
$wsDeDe.Range("H5").NumberFormat = $dateFormat
$wsDeDe.Range("I5").Value = "'"
$wsDeDe.Range("J5").Value = "'"
$wsDeDe.Range("K5").Value = "0001-01-01 00:00:00"
$wsDeDe.Range("K5").NumberFormat = $dateFormat
$wsDeDe.Range("L5").Value = "'"
$wsDeDe.Range("M5").Value = "'True"
$wsDeDe.Range("N5").Value = "'"
$wsDeDe.Range("O5").Value = "'False"
$wsDeDe.Range("P5").Value = "'"

$wsDeDe.Hyperlinks.Add($wsDeDe.Range("A4"), "https://github.com/OpenLocalizationTestOrg/ol-test0-dede/blob/81bd31b2a832704f2ef8133ffad3685011d30cb5/e2e/646596ff-78c4-4b47-b2ef-2067e79379b0.md", "", "", "646596ff-78c4-4b47-b2ef-2067e79379b0.md") | Out-Null
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("A5"), "https://github.com/OpenLocalizationTestOrg/ol-test0-dede/blob/81bd31b2a832704f2ef8133ffad3685011d30cb5/e2e/68bcec46-6cb0-401a-b119-925d1709bc41.md", "", "", "68bcec46-6cb0-401a-b119-925d1709bc41.md") | Out-Null
